$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$updates = @{
    2  = -4
    4  = 8
    7  = 12
    9  = -2
    11 = -1
    15 = 4
    17 = -6
    18 = 2
    22 = -4
    23 = -6
    25 = -2
    27 = -2
    28 = 1
    29 = -3
    30 = -4
    31 = -3
    37 = -6
    41 = -4
    45 = 0
    49 = -10
    51 = -3
    52 = 0
    54 = 0
    56 = 6
    58 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
